# Applies the cryptos-list price/volume refresh described in the commit:
# "Updated cryptos list on Sun Jun 23 16:36:42 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.104.49'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '3.477.90'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'584.62"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = "'131.44"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = "'0.482"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('D9').Value = "'7.62"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.66%  '
$ws.Range('E10').Value = '  -1.72%  '
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '4.068.17'
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').Value = '3.478.34'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('D16').Value = '64.101.80'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = "'24.30"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.72%  '
$ws.Range('D18').Value = "'9.94"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').Value = "'13.42"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.30%  '
$ws.Range('D21').Value = "'384.39"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.56%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').Value = '3.616.60'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('D24').Value = "'74.69"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').Value = "'5.64"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('E27').Value = '  -2.23%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D30').Value = "'7.13"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('E31').Value = '  -4.24%  '
$ws.Range('E32').Value = '  -4.31%  '
$ws.Range('D33').Value = '3.506.22'
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').Value = "'0.152"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('D36').Value = "'22.95"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').Value = "'6.76"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.29%  '
$ws.Range('E39').Value = '  -3.95%  '
$ws.Range('D40').Value = "'162.58"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').Value = "'0.0775"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.21%  '
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').Value = "'41.40"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('D45').Value = "'4.30"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('E46').Value = '  -1.96%  '
$ws.Range('E47').Value = '  -6.52%  '
$ws.Range('E48').Value = '  -3.87%  '
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('D50').Value = "'0.902"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').Value = '2.351.16'
$ws.Range('E51').Value = '  -4.87%  '
